# Applies the periodic cryptos-list refresh (prices / 1h volume %, plus a few
# rank swaps) to Sheet1, mirroring the GitHub Actions data-refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing it to stay a TEXT cell (the
# sheet stores Price as text, e.g. "317.70" / "0.0000137" / "8.59" -- left to
# its own devices Excel would coerce these plain-looking numeric strings into
# real numbers and mangle them (drop trailing zeros, use exponent notation, ...).
# Temporarily marking the cell as Text (@) for the assignment keeps the literal
# string, and restoring the original .Style afterwards leaves formatting untouched.
function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

# Row 2: Bitcoin
$ws.Range("D2").Value = '63.248.46'
$ws.Range("E2").Value = '  +10.68%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '3.478.76'
$ws.Range("E3").Value = '  +6.81%  '

# Row 4: TetherUSD
$ws.Range("E4").Value = '  +0.06%  '

# Row 5: BNB
Set-TextValue $ws.Range("D5") '416.19'
$ws.Range("E5").Value = '  +4.69%  '

# Row 6: Solana
Set-TextValue $ws.Range("D6") '124.34'
$ws.Range("E6").Value = '  +14.22%  '

# Row 7: LidoStakedEther
$ws.Range("D7").Value = '3.474.02'
$ws.Range("E7").Value = '  +6.77%  '

# Row 8: XRP
Set-TextValue $ws.Range("D8") '0.602'
$ws.Range("E8").Value = '  +4.38%  '

# Row 9: USDC
$ws.Range("E9").Value = '  +0.07%  '

# Row 10: Cardano
Set-TextValue $ws.Range("D10") '0.675'
$ws.Range("E10").Value = '  +9.11%  '

# Row 11: Dogecoin
Set-TextValue $ws.Range("D11") '0.131'
$ws.Range("E11").Value = '  +37.35%  '

# Row 12: Avalanche
Set-TextValue $ws.Range("D12") '41.58'
$ws.Range("E12").Value = '  +5.94%  '

# Row 13: TRON
$ws.Range("E13").Value = '  -0.16%  '

# Row 14: WrappedliquidstakedEther2.0
$ws.Range("D14").Value = '4.028.99'
$ws.Range("E14").Value = '  +6.77%  '

# Row 15: Polkadot
Set-TextValue $ws.Range("D15") '8.59'
$ws.Range("E15").Value = '  +4.39%  '

# Row 16: Chainlink
Set-TextValue $ws.Range("D16") '19.95'
$ws.Range("E16").Value = '  +5.29%  '

# Row 17: WrappedEther
$ws.Range("D17").Value = '3.470.33'
$ws.Range("E17").Value = '  +6.70%  '

# Row 18: WrappedBTC
$ws.Range("D18").Value = '63.139.75'
$ws.Range("E18").Value = '  +10.90%  '

# Row 19: Polygon
$ws.Range("E19").Value = '  +0.69%  '

# Row 20: Uniswap
Set-TextValue $ws.Range("D20") '11.13'
$ws.Range("E20").Value = '  +1.24%  '

# Row 21: ShibaInu
Set-TextValue $ws.Range("D21") '0.0000137'
$ws.Range("E21").Value = '  +27.36%  '

# Row 22: ImmutableX
$ws.Range("E22").Value = '  +1.56%  '

# Row 23: BitcoinCash/Litecoin (swap)
$ws.Range("B23").Value = 'Litecoin'
$ws.Range("C23").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue $ws.Range("D23") '83.15'
$ws.Range("E23").Value = '  +12.50%  '

# Row 24: Litecoin/BitcoinCash (swap)
$ws.Range("B24").Value = 'BitcoinCash'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue $ws.Range("D24") '317.70'
$ws.Range("E24").Value = '  +7.74%  '

# Row 25: InternetComputer(DFINITY)
Set-TextValue $ws.Range("D25") '12.97'
$ws.Range("E25").Value = '  +0.35%  '

# Row 26: PancakeSwap
$ws.Range("E26").Value = '  +0.60%  '

# Row 27: EthereumClassic
Set-TextValue $ws.Range("D27") '31.10'
$ws.Range("E27").Value = '  +10.83%  '

# Row 28: RenderToken
Set-TextValue $ws.Range("D28") '7.81'
$ws.Range("E28").Value = '  +5.30%  '

# Row 29: Filecoin
Set-TextValue $ws.Range("D29") '7.91'
$ws.Range("E29").Value = '  +0.29%  '

# Row 30: Kaspa
Set-TextValue $ws.Range("D30") '0.176'
$ws.Range("E30").Value = '  +4.26%  '

# Row 31: LEO
$ws.Range("E31").Value = '  -1.51%  '

# Row 32: Hedera
$ws.Range("E32").Value = '  +3.12%  '

# Row 33: Cosmos
Set-TextValue $ws.Range("D33") '11.65'
$ws.Range("E33").Value = '  +4.21%  '

# Row 34: InjectiveProtocol/Toncoin (swap)
$ws.Range("B34").Value = 'Toncoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws.Range("D34") '2.57'
$ws.Range("E34").Value = '  +19.94%  '

# Row 35: Toncoin/InjectiveProtocol (swap)
$ws.Range("B35").Value = 'InjectiveProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range("D35") '42.04'
$ws.Range("E35").Value = '  +3.37%  '

# Row 36: Dai
$ws.Range("E36").Value = '  +0.58%  '

# Row 37: VeChain
$ws.Range("E37").Value = '  +0.29%  '

# Row 38: OKB
Set-TextValue $ws.Range("D38") '52.24'
$ws.Range("E38").Value = '  +1.89%  '

# Row 39: LidoDAOToken/FirstDigitalUSD (swap)
$ws.Range("B39").Value = 'FirstDigitalUSD'
$ws.Range("C39").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws.Range("D39") '0.999'
$ws.Range("E39").Value = '  -0.09%  '

# Row 40: FirstDigitalUSD/LidoDAOToken (swap)
$ws.Range("B40").Value = 'LidoDAOToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue $ws.Range("D40") '3.48'
$ws.Range("E40").Value = '  +0.99%  '

# Row 41: Stacks
Set-TextValue $ws.Range("D41") '3.05'
$ws.Range("E41").Value = '  +1.98%  '

# Row 42: ARBITRUM
Set-TextValue $ws.Range("D42") '2.01'
$ws.Range("E42").Value = '  +7.76%  '

# Row 43: Stellar
$ws.Range("E43").Value = '  +4.81%  '

# Row 44: Monero
Set-TextValue $ws.Range("D44") '135.11'
$ws.Range("E44").Value = '  -1.24%  '

# Row 45: Celestia
Set-TextValue $ws.Range("D45") '17.29'
$ws.Range("E45").Value = '  +3.49%  '

# Row 46: TheGraph
Set-TextValue $ws.Range("D46") '0.285'
$ws.Range("E46").Value = '  +0.94%  '

# Row 47: NEARProtocol
Set-TextValue $ws.Range("D47") '3.93'
$ws.Range("E47").Value = '  +0.65%  '

# Row 48: WEMIXToken
Set-TextValue $ws.Range("D48") '2.25'
$ws.Range("E48").Value = '  +2.12%  '

# Row 49: EnergySwap
Set-TextValue $ws.Range("D49") '22.21'
$ws.Range("E49").Value = '  -0.55%  '

# Row 50: RocketPoolETH
$ws.Range("D50").Value = '3.820.20'
$ws.Range("E50").Value = '  +6.70%  '

# Row 51: Maker
$ws.Range("D51").Value = '2.185.98'
$ws.Range("E51").Value = '  +1.97%  '
